# ADD results from server
# Update computed result values on row 2 of each yearly sheet with
# refreshed figures received from the server.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param(
        [string]$SheetName,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $addr = "$($col)2"
        $ws.Range($addr).Value = $Values[$col]
    }
}

Set-RowValues "2025" @{
    "N" = 7155.075790473336
    "O" = 6980.325566461754
}

Set-RowValues "2030" @{
    "B" = 5707.815717280662
    "I" = 44492.05901988943
    "L" = 66334.06707325629
    "M" = 21991.42050229464
    "O" = 12079.40905079305
}

Set-RowValues "2035" @{
    "A" = 2927.360317916481
    "B" = 7940.887964949257
    "E" = 67179.99183625776
    "I" = 59530.75343380851
    "L" = 66334.06707325629
    "M" = 25547.11936466757
    "N" = 15117.91059331085
    "O" = 14761.05415301146
}

Set-RowValues "2040" @{
    "A" = 2927.360317916481
    "B" = 7940.887964949257
    "E" = 67179.99183625776
    "I" = 59530.75343380851
    "L" = 66334.06707325629
    "M" = 25547.11936466757
    "N" = 15222.78766604848
    "O" = 14761.05415301146
}

Set-RowValues "2045" @{
    "A" = 6352.985609279765
    "B" = 7940.887964949257
    "E" = 67179.99183625776
    "I" = 59530.75343380851
    "L" = 66334.06707325629
    "M" = 25547.11936466757
    "N" = 15767.51521749871
    "O" = 17096.52013936021
}

Set-RowValues "2050" @{
    "A" = 6352.985609279765
    "B" = 7940.887964949257
    "E" = 67179.99183625776
    "I" = 59530.75343380851
    "L" = 66334.06707325629
    "M" = 25547.11936466757
    "N" = 15767.51521749871
    "O" = 17096.52013936021
}
